# Refresh the coin price table (cols B-D: name/link/price, col E: 1h volume %)
# with the latest scraped figures (GitHub Actions cron update).
#
# Price values (col D) are sometimes plain decimals ("0.9994", "1.000", ...).
# Assigning those straight to Range.Value lets Excel auto-detect them as
# numbers, which silently drops formatting such as trailing zeros. Since the
# sheet stores every Price/Volume cell as text, Set-TextValue below detects
# plain-numeric strings and writes them with a leading apostrophe (the normal
# Excel "store number as text" idiom) so they round-trip byte-for-byte.
function Set-TextValue($range, $value) {
    if ($value -match '^-?\d+(\.\d+)?$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") '29.116.96'

Set-TextValue $ws.Range("D3") '1.835.86'
$ws.Range("E3").Value = '  -1.43%  '

Set-TextValue $ws.Range("D4") '0.9994'
$ws.Range("E4").Value = '  -0.03%  '

Set-TextValue $ws.Range("D5") '239.73'
$ws.Range("E5").Value = '  -2.40%  '

Set-TextValue $ws.Range("D7") '1.000'
$ws.Range("E7").Value = '  +0.00%  '

Set-TextValue $ws.Range("D8") '0.2987'

Set-TextValue $ws.Range("D9") '0.07444'
$ws.Range("E9").Value = '  -3.84%  '

Set-TextValue $ws.Range("D10") '23.13'
$ws.Range("E10").Value = '  -2.39%  '

$ws.Range("E11").Value = '  -1.60%  '

Set-TextValue $ws.Range("D12") '1.822.64'
$ws.Range("E12").Value = '  -2.07%  '

Set-TextValue $ws.Range("D13") '5.023'
$ws.Range("E13").Value = '  -2.67%  '

Set-TextValue $ws.Range("D14") '0.6784'
$ws.Range("E14").Value = '  -2.07%  '

Set-TextValue $ws.Range("D15") '86.85'
$ws.Range("E15").Value = '  -5.93%  '

Set-TextValue $ws.Range("D16") '6.155'
$ws.Range("E16").Value = '  -6.30%  '

Set-TextValue $ws.Range("D17") '29.111.80'
$ws.Range("E17").Value = '  -1.72%  '

Set-TextValue $ws.Range("D18") '0.000008230'
$ws.Range("E18").Value = '  -1.69%  '

Set-TextValue $ws.Range("D19") '229.06'
$ws.Range("E19").Value = '  -5.36%  '

$ws.Range("E20").Value = '  -2.24%  '

Set-TextValue $ws.Range("D21") '0.9993'
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("E22").Value = '  -3.73%  '

Set-TextValue $ws.Range("D23") '1.000'

Set-TextValue $ws.Range("D24") '161.30'
$ws.Range("E24").Value = '  +1.04%  '

Set-TextValue $ws.Range("D25") '0.1431'
$ws.Range("E25").Value = '  -5.19%  '

$ws.Range("E26").Value = '  -2.37%  '

$ws.Range("E27").Value = '  -1.58%  '

Set-TextValue $ws.Range("D28") '1.497'
$ws.Range("E28").Value = '  -2.46%  '

$ws.Range("E29").Value = '  -0.36%  '

Set-TextValue $ws.Range("D30") '4.140'
$ws.Range("E30").Value = '  -1.40%  '

$ws.Range("E31").Value = '  -0.46%  '

$ws.Range("E32").Value = '  +4.41%  '

$ws.Range("E33").Value = '  -3.90%  '

$ws.Range("E34").Value = '  -3.04%  '

Set-TextValue $ws.Range("D35") '1.131'
$ws.Range("E35").Value = '  -2.41%  '

Set-TextValue $ws.Range("D36") '2.683'
$ws.Range("E36").Value = '  -0.12%  '

Set-TextValue $ws.Range("D37") '1.313.28'
$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("E38").Value = '  -3.22%  '

$ws.Range("E39").Value = '  -0.91%  '

Set-TextValue $ws.Range("D40") '0.9374'
$ws.Range("E40").Value = '  -2.75%  '

Set-TextValue $ws.Range("D41") '6.062'
$ws.Range("E41").Value = '  +2.02%  '

$ws.Range("B42").Value = 'XinFinNetwork'
$ws.Range("C42").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
Set-TextValue $ws.Range("D42") '0.08450'
$ws.Range("E42").Value = '  +33.30%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D43") '104.93'
$ws.Range("E43").Value = '  -1.50%  '

$ws.Range("E44").Value = '  -0.09%  '

Set-TextValue $ws.Range("D45") '1.981.44'
$ws.Range("E45").Value = '  -1.31%  '

Set-TextValue $ws.Range("D46") '0.5178'
$ws.Range("E46").Value = '  -0.69%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D47") '1.768'
$ws.Range("E47").Value = '  -1.27%  '

Set-TextValue $ws.Range("D48") '63.92'
$ws.Range("E48").Value = '  -1.19%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D49") '9.372'
$ws.Range("E49").Value = '  -4.18%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.05927'
$ws.Range("E50").Value = '  +0.25%  '

$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D51") '6.890'
$ws.Range("E51").Value = '  -1.69%  '
